$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Prix Spot" -> add column AA (10-jul) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("Z1").Copy($ws1.Range("AA1"))
$ws1.Range("AA1").Value = "10-jul"

$aaValues = @(100.91, 95.73, 89.09999999999999, 66.89, 48.68, 60.43, 86.3, 90.47, 91.42, 79.29000000000001, 33.45, 19.99, 40.5, 17.59, 11.9, 16.49, 39.27, 49.39, 70.20999999999999, 99.98999999999999, 110.04, 100.16, 111, 93.65000000000001)

for ($i = 0; $i -lt $aaValues.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 27).Value = $aaValues[$i]
}

# --- Sheet 2: "Gaz" -> add row 24 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Cells.Item(24, 1).NumberFormat = "@"
$ws2.Cells.Item(24, 1).Value = "2025-07-08"
$ws2.Cells.Item(24, 1).Style = "Normal"
$ws2.Cells.Item(24, 2).Value = 33.85

# --- Sheet 3: "CO2" -> add row 24 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Cells.Item(24, 1).NumberFormat = "@"
$ws3.Cells.Item(24, 1).Value = "2025-07-08"
$ws3.Cells.Item(24, 1).Style = "Normal"
$ws3.Cells.Item(24, 2).Value = 69.7
